# Daily attendance processing - 2025-11-13 15:47:52
# Re-sort the "Recorded By" (column G) comma-separated name lists into
# ordinal (ASCII) ascending order, e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

function Compare-Ordinal($left, $right) {
    $lenLeft = $left.Length
    $lenRight = $right.Length
    $minLen = [Math]::Min($lenLeft, $lenRight)
    for ($pos = 0; $pos -lt $minLen; $pos++) {
        $codeLeft = [int][char]$left[$pos]
        $codeRight = [int][char]$right[$pos]
        if ($codeLeft -lt $codeRight) { return -1 }
        if ($codeLeft -gt $codeRight) { return 1 }
    }
    if ($lenLeft -lt $lenRight) { return -1 }
    if ($lenLeft -gt $lenRight) { return 1 }
    return 0
}

function Sort-Ordinal($items) {
    $arrS = @($items)
    $cntS = $arrS.Count
    for ($idxI = 1; $idxI -lt $cntS; $idxI++) {
        $keyS = $arrS[$idxI]
        $idxJ = $idxI - 1
        while ($idxJ -ge 0 -and (Compare-Ordinal $arrS[$idxJ] $keyS) -gt 0) {
            $arrS[$idxJ + 1] = $arrS[$idxJ]
            $idxJ = $idxJ - 1
        }
        $arrS[$idxJ + 1] = $keyS
    }
    return $arrS
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($rowNum = 2; $rowNum -le $lastRow; $rowNum++) {
    $cell = $ws.Cells.Item($rowNum, 7)
    $raw = $cell.Value()
    if ($raw -eq $null) { continue }
    if ($raw -eq "") { continue }

    $rawParts = $raw -split ","
    $trimmedParts = @()
    foreach ($part in $rawParts) {
        $trimmedParts += $part.Trim()
    }

    $sortedParts = Sort-Ordinal $trimmedParts
    $newValue = $sortedParts -join ", "

    if ($newValue -ne $raw) {
        $cell.Value = $newValue
    }
}
